# Auto-generated edit script: refresh cryptos Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$style = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.213.25'
$ws.Range("D2").Style = $style
$ws.Range("E2").Value = '  -0.26%  '

$style = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.643.30'
$ws.Range("D3").Style = $style
$ws.Range("E3").Value = '  -0.13%  '

$ws.Range("E4").Value = '  +0.01%  '

$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '595.60'
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = '  -0.44%  '

$style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '158.74'
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = '  +2.44%  '

$ws.Range("E7").Value = '  -0.01%  '

$style = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.541'
$ws.Range("D8").Style = $style
$ws.Range("E8").Value = '  -1.02%  '

$ws.Range("E9").Value = '  -2.76%  '

$style = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.156'
$ws.Range("D10").Style = $style
$ws.Range("E10").Value = '  -1.16%  '

$ws.Range("E11").Value = '  +0.15%  '

$ws.Range("E12").Value = '  -1.16%  '

$style = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '27.89'
$ws.Range("D13").Style = $style
$ws.Range("E13").Value = '  -1.30%  '

$style = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.127.62'
$ws.Range("D14").Style = $style
$ws.Range("E14").Value = '  +0.03%  '

$ws.Range("E15").Value = '  -3.65%  '

$style = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '68.080.90'
$ws.Range("D16").Style = $style
$ws.Range("E16").Value = '  -0.38%  '

$style = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.638.43'
$ws.Range("D17").Style = $style
$ws.Range("E17").Value = '  -0.36%  '

$style = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.34'
$ws.Range("D18").Style = $style
$ws.Range("E18").Value = '  -0.96%  '

$style = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '358.76'
$ws.Range("D19").Style = $style
$ws.Range("E19").Value = '  -1.74%  '

$style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.39'
$ws.Range("D20").Style = $style
$ws.Range("E20").Value = '  -1.91%  '

$ws.Range("E21").Value = '  +0.33%  '

$ws.Range("E22").Value = '  -3.48%  '

$ws.Range("E23").Value = '  -0.70%  '

$style = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '74.50'
$ws.Range("D24").Style = $style
$ws.Range("E24").Value = '  -0.08%  '

$ws.Range("E25").Value = '  +0.04%  '

$style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.70'
$ws.Range("D26").Style = $style
$ws.Range("E26").Value = '  -1.13%  '

$style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.777.72'
$ws.Range("D27").Style = $style
$ws.Range("E27").Value = '  +0.08%  '

$ws.Range("E28").Value = '  -3.35%  '

$ws.Range("E29").Value = '  -0.02%  '

$style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '559.31'
$ws.Range("D30").Style = $style
$ws.Range("E30").Value = '  -2.46%  '

$style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.98'
$ws.Range("D31").Style = $style
$ws.Range("E31").Value = '  -2.64%  '

$ws.Range("E32").Value = '  -2.77%  '

$style = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.87'
$ws.Range("D33").Style = $style
$ws.Range("E33").Value = '  +0.60%  '

$ws.Range("E34").Value = '  +1.77%  '

$ws.Range("E35").Value = '  -0.02%  '

$ws.Range("E36").Value = '  -3.25%  '

$style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '159.23'
$ws.Range("D37").Style = $style
$ws.Range("E37").Value = '  -1.06%  '

$style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '19.65'
$ws.Range("D38").Style = $style
$ws.Range("E38").Value = '  +1.13%  '

$style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.369'
$ws.Range("D39").Style = $style
$ws.Range("E39").Value = '  -1.42%  '

$style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.85'
$ws.Range("D40").Style = $style
$ws.Range("E40").Value = '  -2.86%  '

$ws.Range("E41").Value = '  -2.35%  '

$style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.59'
$ws.Range("D42").Style = $style
$ws.Range("E42").Value = '  -2.73%  '

$style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0₆0318'
$ws.Range("D43").Style = $style
$ws.Range("E43").Value = '  -5.85%  '

$ws.Range("E44").Value = '  +0.03%  '

$style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '157.26'
$ws.Range("D45").Style = $style
$ws.Range("E45").Value = '  +0.15%  '

$style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.78'
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = '  +0.10%  '

$ws.Range("E47").Value = '  -0.36%  '

$style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.68'
$ws.Range("D48").Style = $style
$ws.Range("E48").Value = '  -2.35%  '

$ws.Range("E49").Value = '  -2.12%  '

$style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.573'
$ws.Range("D50").Style = $style
$ws.Range("E50").Value = '  +0.54%  '

$style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.614'
$ws.Range("D51").Style = $style
$ws.Range("E51").Value = '  -0.46%  '

